# v2.2.1.9 - New command to set action speed
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "V2 Command" sheet: insert a new row 48 for the "Set Action Speed"
#    command, fix D40, and update the view state.
# ---------------------------------------------------------------------------
$cmd = $wb.Worksheets.Item("V2 Command")

# Insert a new row at 48 - everything below (old rows 48-71) shifts to 49-72.
$cmd.Rows.Item(48).Insert()

# New row content.
$cmd.Range("B48").Value = 43
$cmd.Range("C48").Value = "Set Action Speed"
$cmd.Range("D48").Value = "Yes {3}"
$cmd.Range("E48").Value = "Speed (time = time * 100 / speed)"
$cmd.Range("F48").Value = "A9 9A 03 43 64 AA ED"
$cmd.Range("G48").Value = "A9 9A 03 43 C8 0E ED"

# Data correction on existing row 40 ("Set Head LED").
$cmd.Range("D40").Value = "Yes {3}"

# ---------------------------------------------------------------------------
# 2) Activate the "V2 Command" sheet (it becomes the saved active tab).
# ---------------------------------------------------------------------------
$cmd.Activate()
$cmd.Application.ActiveWindow.ScrollRow = 20
$cmd.Range("G48").Select()

Write-Host "done"
